# Refresh the cryptos price table (Price / Volume(1h) columns, plus two
# coin-row swaps) to match the latest scrape.
#
# Note: several "Price" values (column D) are plain decimal-looking text
# (e.g. "316.90", "0.0844") that must stay TEXT, not become numbers -
# otherwise trailing zeros are lost (e.g. "316.90" -> 316.9) and Excel's
# own float rendering creeps in (e.g. "92.72" -> 92.719999999999999).
# Assigning with a leading apostrophe forces text storage (quote-prefix),
# then resetting .Style back to 'Normal' strips the visual quote-prefix
# formatting again so the cell ends up styled exactly like its neighbours.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.748.30'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range('D3').Value = '2.467.95'
$ws.Range('E3').Value = '  -0.79%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = "'316.90"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('D6').Value = "'92.72"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('E7').Value = '  +0.80%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  +3.45%  '
$ws.Range('D10').Value = "'32.79"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.37%  '
$ws.Range('D11').Value = "'0.0844"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.96%  '
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('D13').Value = '2.850.43'
$ws.Range('E13').Value = '  -0.51%  '
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('D15').Value = "'15.77"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.17%  '
$ws.Range('D16').Value = '2.497.79'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').Value = "'0.780"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.59%  '
$ws.Range('D18').Value = '41.744.96'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').Value = "'6.50"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.71%  '
$ws.Range('D20').Value = '0.0₃0955'
$ws.Range('E20').Value = '  +3.53%  '
$ws.Range('D21').Value = "'11.67"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.49%  '
$ws.Range('D22').Value = "'71.15"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.76%  '
$ws.Range('D23').Value = "'239.30"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.59%  '
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('E25').Value = '  +1.24%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('D28').Value = "'2.26"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.15%  '
$ws.Range('D29').Value = "'9.79"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.46%  '
$ws.Range('D30').Value = "'36.01"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.79%  '
$ws.Range('D31').Value = "'155.92"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range('D32').Value = "'5.51"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.80%  '
$ws.Range('E33').Value = '  +0.19%  '
$ws.Range('E34').Value = '  +1.50%  '
$ws.Range('E35').Value = '  +1.88%  '
$ws.Range('D36').Value = "'17.62"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.41%  '
$ws.Range('D37').Value = "'2.90"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.72%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = "'0.104"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').Value = "'0.115"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.11%  '
$ws.Range('D40').Value = "'1.80"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.81%  '
$ws.Range('D41').Value = "'4.03"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.95%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').Value = '1.976.71'
$ws.Range('E43').Value = '  +0.70%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = "'19.05"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.85%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = "'0.0284"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('E46').Value = '  -0.33%  '
$ws.Range('E47').Value = '  +2.04%  '
$ws.Range('D48').Value = '2.703.24'
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('D49').Value = "'96.96"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('D50').Value = "'67.24"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('D51').Value = "'73.31"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.10%  '
